$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("E2").Value = 299457.1989099024
$ws.Range("G2").Value = 80959.25712661653
$ws.Range("I2").Value = 136981.2725754988
$ws.Range("L2").Value = 540337.8794259601
$ws.Range("M2").Value = 105561.3095757
$ws.Range("N2").Value = 75634.14130920739
$ws.Range("O2").Value = 75817.64920437815

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 42918.32072666658
$ws.Range("E2").Value = 281775.2215551065
$ws.Range("I2").Value = 134723.1329735118
$ws.Range("L2").Value = 410592.6555366496
$ws.Range("M2").Value = 100773.5540459725
$ws.Range("N2").Value = 32503.43204176008
$ws.Range("O2").Value = 20101.74547350489

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 15243.27458338266
$ws.Range("B2").Value = 28182.92772446093
$ws.Range("E2").Value = 138653.2715386877
$ws.Range("I2").Value = 172494.5760716226
$ws.Range("M2").Value = 30391.66382300033
$ws.Range("N2").Value = 49794.70591665693
$ws.Range("O2").Value = 55798.41569072519
